# Update (Analyze PO & Forecast)
# Applies refreshed forecast numbers to the "Forecast Comparison" sheet
# (MyForecast column D) and the corresponding roll-up figures on the
# "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- "Forecast Comparison": MyForecast (column D) updates ---
$wsForecast.Range("D2").Value  = 18   # W1
$wsForecast.Range("D3").Value  = 18   # W2
$wsForecast.Range("D7").Value  = 17   # W6
$wsForecast.Range("D8").Value  = 17   # W7
$wsForecast.Range("D9").Value  = 16   # W8
$wsForecast.Range("D10").Value = 15   # W9
$wsForecast.Range("D11").Value = 13   # W10
$wsForecast.Range("D12").Value = 11   # W11
$wsForecast.Range("D14").Value = 11   # W13
$wsForecast.Range("D15").Value = 13   # W14
$wsForecast.Range("D16").Value = 12   # W15
$wsForecast.Range("D17").Value = 11   # W16

# --- "Summary": updated roll-up metrics ---
# These "Value" column cells are stored as text in the workbook, so a
# leading apostrophe is used to keep numeric-looking / date-looking
# entries as text instead of letting them be auto-converted to a number
# or date serial.
$wsSummary.Range("B9").Value  = "'234"        # Total Forecast (16 Weeks)
$wsSummary.Range("B10").Value = "'138"        # Total Forecast (8 Weeks)
$wsSummary.Range("B11").Value = "'70"         # Total Forecast (4 Weeks)
$wsSummary.Range("B14").Value = "'10"         # Min Forecast
$wsSummary.Range("B15").Value = "'2025-04-13" # Min Forecast Week
